$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column header: clarify Role labels ---
$ws.Range("A1").Value = "Role (PI/Sub I)"

# --- Row 2 (Principal investigator record) ---
$ws.Range("A2").Value = "PI"
$ws.Range("K2").Value = "avenue road"
$ws.Range("L2").ClearContents()
$ws.Range("O2").Value = 889245
$ws.Range("P2").Value = "USA"

# --- Sub-investigator rows: "sub" -> "sub I" ---
$ws.Range("A3").Value = "sub I"
$ws.Range("A4").Value = "sub I"
$ws.Range("A5").Value = "sub I"
$ws.Range("A7").Value = "sub I"
$ws.Range("A8").Value = "sub I"
$ws.Range("A9").Value = "sub I"

# --- Principal investigator rows: "Principal" -> "PI" ---
$ws.Range("A6").Value = "PI"
$ws.Range("A10").Value = "PI"

# --- Scroll/selection state as left by the editor ---
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("O3").Select()
